$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newA = @(
45969,45969.01041666666,45969.02083333334,45969.03125,45969.04166666666,45969.05208333334,45969.0625,45969.07291666666,45969.08333333334,45969.09375,45969.10416666666,45969.11458333334,45969.125,45969.13541666666,45969.14583333334,45969.15625,45969.16666666666,45969.17708333334,45969.1875,45969.19791666666,45969.20833333334,45969.21875,45969.22916666666,45969.23958333334,45969.25,45969.26041666666,45969.27083333334,45969.28125,45969.29166666666,45969.30208333334,45969.3125,45969.32291666666,45969.33333333334,45969.34375,45969.35416666666,45969.36458333334,45969.375,45969.38541666666,45969.39583333334,45969.40625,45969.41666666666,45969.42708333334,45969.4375,45969.44791666666,45969.45833333334,45969.46875,45969.47916666666,45969.48958333334,45969.5,45969.51041666666,45969.52083333334,45969.53125,45969.54166666666,45969.55208333334,45969.5625,45969.57291666666,45969.58333333334,45969.59375,45969.60416666666,45969.61458333334,45969.625,45969.63541666666,45969.64583333334,45969.65625,45969.66666666666,45969.67708333334,45969.6875,45969.69791666666,45969.70833333334,45969.71875,45969.72916666666,45969.73958333334,45969.75,45969.76041666666,45969.77083333334,45969.78125,45969.79166666666,45969.80208333334,45969.8125,45969.82291666666,45969.83333333334,45969.84375,45969.85416666666,45969.86458333334,45969.875,45969.88541666666,45969.89583333334,45969.90625,45969.91666666666,45969.92708333334,45969.9375,45969.94791666666,45969.95833333334,45969.96875,45969.97916666666,45969.98958333334,45970,45970.01041666666,45970.02083333334,45970.03125,45970.04166666666,45970.05208333334,45970.0625,45970.07291666666,45970.08333333334,45970.09375,45970.10416666666,45970.11458333334,45970.125,45970.13541666666,45970.14583333334,45970.15625,45970.16666666666,45970.17708333334,45970.1875,45970.19791666666,45970.20833333334,45970.21875,45970.22916666666,45970.23958333334,45970.25,45970.26041666666,45970.27083333334,45970.28125,45970.29166666666,45970.30208333334,45970.3125,45970.32291666666,45970.33333333334,45970.34375,45970.35416666666,45970.36458333334,45970.375,45970.38541666666,45970.39583333334,45970.40625,45970.41666666666,45970.42708333334,45970.4375,45970.44791666666,45970.45833333334,45970.46875,45970.47916666666,45970.48958333334,45970.5,45970.51041666666,45970.52083333334,45970.53125,45970.54166666666,45970.55208333334,45970.5625,45970.57291666666,45970.58333333334,45970.59375,45970.60416666666,45970.61458333334,45970.625,45970.63541666666,45970.64583333334,45970.65625,45970.66666666666,45970.67708333334,45970.6875,45970.69791666666,45970.70833333334,45970.71875,45970.72916666666,45970.73958333334,45970.75,45970.76041666666,45970.77083333334,45970.78125,45970.79166666666,45970.80208333334,45970.8125,45970.82291666666,45970.83333333334,45970.84375,45970.85416666666,45970.86458333334,45970.875,45970.88541666666,45970.89583333334,45970.90625,45970.91666666666,45970.92708333334,45970.9375,45970.94791666666,45970.95833333334,45970.96875,45970.97916666666,45970.98958333334
)

$newB = @(
5600,5569,5557,5475,5469,5449,5412,5365,5403,5364,5399,5391,5423,5402,5421,5401,5457,5475,5481,5526,5485,5534,5588,5636,5744,5821,5843,5915,5985,6024,6101,6070,6135,6216,6162,6193,6169,6216,6234,6306,6141,6195,6223,6282,6210,6222,6267,6266,6253,6218,6257,6248,6286,6255,6282,6283,6309,6305,6380,6418,6511,6533,6611,6703,6806,6874,6954,7006,7002,6964,6913,6930,6912,6842,6818,6791,6752,6713,6668,6630,6516,6440,6318,6230,6092,6015,5918,5795,5712,5607,5517,5439,5424,5455,5430,5337,5216,5126,5129,5037,5048,5060,5022,5029,4966,4994,4913,4884,4913,4897,4904,4886,4956,4906,4948,4986,5046,5055,5116,5102,5093,5143,5207,5258,5269,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0
)

$newD = @(
"08.11.20251","08.11.20252","08.11.20253","08.11.20254","08.11.20255","08.11.20256","08.11.20257","08.11.20258","08.11.20259","08.11.202510","08.11.202511","08.11.202512","08.11.202513","08.11.202514","08.11.202515","08.11.202516","08.11.202517","08.11.202518","08.11.202519","08.11.202520","08.11.202521","08.11.202522","08.11.202523","08.11.202524","08.11.202525","08.11.202526","08.11.202527","08.11.202528","08.11.202529","08.11.202530","08.11.202531","08.11.202532","08.11.202533","08.11.202534","08.11.202535","08.11.202536","08.11.202537","08.11.202538","08.11.202539","08.11.202540","08.11.202541","08.11.202542","08.11.202543","08.11.202544","08.11.202545","08.11.202546","08.11.202547","08.11.202548","08.11.202549","08.11.202550","08.11.202551","08.11.202552","08.11.202553","08.11.202554","08.11.202555","08.11.202556","08.11.202557","08.11.202558","08.11.202559","08.11.202560","08.11.202561","08.11.202562","08.11.202563","08.11.202564","08.11.202565","08.11.202566","08.11.202567","08.11.202568","08.11.202569","08.11.202570","08.11.202571","08.11.202572","08.11.202573","08.11.202574","08.11.202575","08.11.202576","08.11.202577","08.11.202578","08.11.202579","08.11.202580","08.11.202581","08.11.202582","08.11.202583","08.11.202584","08.11.202585","08.11.202586","08.11.202587","08.11.202588","08.11.202589","08.11.202590","08.11.202591","08.11.202592","08.11.202593","08.11.202594","08.11.202595","08.11.202596","09.11.20251","09.11.20252","09.11.20253","09.11.20254","09.11.20255","09.11.20256","09.11.20257","09.11.20258","09.11.20259","09.11.202510","09.11.202511","09.11.202512","09.11.202513","09.11.202514","09.11.202515","09.11.202516","09.11.202517","09.11.202518","09.11.202519","09.11.202520","09.11.202521","09.11.202522","09.11.202523","09.11.202524","09.11.202525","09.11.202526","09.11.202527","09.11.202528","09.11.202529","09.11.202530","09.11.202531","09.11.202532","09.11.202533","09.11.202534","09.11.202535","09.11.202536","09.11.202537","09.11.202538","09.11.202539","09.11.202540","09.11.202541","09.11.202542","09.11.202543","09.11.202544","09.11.202545","09.11.202546","09.11.202547","09.11.202548","09.11.202549","09.11.202550","09.11.202551","09.11.202552","09.11.202553","09.11.202554","09.11.202555","09.11.202556","09.11.202557","09.11.202558","09.11.202559","09.11.202560","09.11.202561","09.11.202562","09.11.202563","09.11.202564","09.11.202565","09.11.202566","09.11.202567","09.11.202568","09.11.202569","09.11.202570","09.11.202571","09.11.202572","09.11.202573","09.11.202574","09.11.202575","09.11.202576","09.11.202577","09.11.202578","09.11.202579","09.11.202580","09.11.202581","09.11.202582","09.11.202583","09.11.202584","09.11.202585","09.11.202586","09.11.202587","09.11.202588","09.11.202589","09.11.202590","09.11.202591","09.11.202592","09.11.202593","09.11.202594","09.11.202595","09.11.202596"
)

for ($i = 0; $i -lt 192; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $newA[$i]
    $ws.Cells.Item($row, 2).Value = $newB[$i]
    $ws.Cells.Item($row, 4).Value = $newD[$i]
}